$d = $word.ActiveDocument

# Find the "Author" styled paragraph that reads exactly "Edison Achalma"
# (the one right after the main title) and add a new "Author" styled
# paragraph right after it with the author's institutional affiliation.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    $styleName = $p.Range.ParagraphStyle.NameLocal
    if ($styleName -eq "Author" -and $text -eq "Edison Achalma") {
        $following = $p.Next()

        # Split a new, empty paragraph right before $following (i.e.
        # immediately after the "Edison Achalma" paragraph) without
        # touching "Edison Achalma" itself.
        $splitPoint = $following.Range.Duplicate
        $splitPoint.Collapse(1)  # wdCollapseStart
        $splitPoint.InsertParagraphAfter()

        $newPara = $p.Next()
        $newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
        $newPara.Style = "Author"
        break
    }
}
